$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'37.080.95"
$ws.Range('E2').Value = "  -0.72%  "
$ws.Range('D3').Value = "'1.993.51"
$ws.Range('E3').Value = "  -1.54%  "
$ws.Range('E4').Value = "  -0.21%  "
$ws.Range('D5').Value = "'264.12"
$ws.Range('E5').Value = "  +7.13%  "
$ws.Range('D6').Value = "'0.608"
$ws.Range('E6').Value = "  -1.88%  "
$ws.Range('E7').Value = "  +0.10%  "
$ws.Range('D8').Value = "'55.78"
$ws.Range('E8').Value = "  -3.84%  "
$ws.Range('D9').Value = "'0.375"
$ws.Range('E9').Value = "  -3.54%  "
$ws.Range('D10').Value = "'0.0761"
$ws.Range('E10').Value = "  -4.87%  "
$ws.Range('E11').Value = "  -3.29%  "
$ws.Range('B12').Value = "Chainlink"
$ws.Range('C12').Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range('D12').Value = "'14.14"
$ws.Range('E12').Value = "  -5.25%  "
$ws.Range('B13').Value = "WrappedliquidstakedEther2.0"
$ws.Range('C13').Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range('D13').Value = "'2.281.09"
$ws.Range('E13').Value = "  -1.90%  "
$ws.Range('D14').Value = "'22.03"
$ws.Range('E14').Value = "  +2.28%  "
$ws.Range('D15').Value = "'0.769"
$ws.Range('E15').Value = "  -7.74%  "
$ws.Range('D16').Value = "'5.17"
$ws.Range('E16').Value = "  -3.97%  "
$ws.Range('D17').Value = "'1.980.32"
$ws.Range('E17').Value = "  -2.65%  "
$ws.Range('D18').Value = "'36.876.89"
$ws.Range('E18').Value = "  -1.17%  "
$ws.Range('D19').Value = "'69.66"
$ws.Range('E19').Value = "  -0.47%  "
$ws.Range('D20').Value = "'0.0₃0825"
$ws.Range('E20').Value = "  -3.30%  "
$ws.Range('D21').Value = "'234.08"
$ws.Range('E21').Value = "  +2.68%  "
$ws.Range('D22').Value = "'5.06"
$ws.Range('E22').Value = "  -2.71%  "
$ws.Range('E23').Value = "  -0.01%  "
$ws.Range('D24').Value = "'2.59"
$ws.Range('E24').Value = "  +2.28%  "
$ws.Range('E25').Value = "  +1.10%  "
$ws.Range('D26').Value = "'165.56"
$ws.Range('E26').Value = "  +1.29%  "
$ws.Range('D27').Value = "'8.82"
$ws.Range('E27').Value = "  -3.76%  "
$ws.Range('D28').Value = "'19.32"
$ws.Range('E28').Value = "  -2.52%  "
$ws.Range('E29').Value = "  -7.74%  "
$ws.Range('D30').Value = "'1.31"
$ws.Range('E30').Value = "  -3.28%  "
$ws.Range('E31').Value = "  -1.97%  "
$ws.Range('D32').Value = "'4.52"
$ws.Range('E32').Value = "  -4.85%  "
$ws.Range('D33').Value = "'0.0620"
$ws.Range('E33').Value = "  -6.99%  "
$ws.Range('D34').Value = "'4.34"
$ws.Range('E34').Value = "  -4.89%  "
$ws.Range('D35').Value = "'2.40"
$ws.Range('E35').Value = "  -2.72%  "
$ws.Range('D36').Value = "'3.50"
$ws.Range('E36').Value = "  -1.00%  "
$ws.Range('E37').Value = "  -0.62%  "
$ws.Range('E38').Value = "  -0.19%  "
$ws.Range('D39').Value = "'5.34"
$ws.Range('E39').Value = "  -0.77%  "
$ws.Range('E40').Value = "  +2.06%  "
$ws.Range('E41').Value = "  +0.38%  "
$ws.Range('D42').Value = "'1.435.61"
$ws.Range('E42').Value = "  +2.71%  "
$ws.Range('D43').Value = "'0.0911"
$ws.Range('E43').Value = "  -6.01%  "
$ws.Range('D44').Value = "'0.0207"
$ws.Range('E44').Value = "  -4.78%  "
$ws.Range('D45').Value = "'89.08"
$ws.Range('E45').Value = "  -1.83%  "
$ws.Range('D46').Value = "'15.41"
$ws.Range('E46').Value = "  -5.93%  "
$ws.Range('E47').Value = "  -2.25%  "
$ws.Range('E48').Value = "  +1.26%  "
$ws.Range('D49').Value = "'6.79"
$ws.Range('E49').Value = "  -8.92%  "
$ws.Range('D50').Value = "'2.173.85"
$ws.Range('E50').Value = "  -2.02%  "
$ws.Range('D51').Value = "'1.92"
$ws.Range('E51').Value = "  -8.01%  "
